$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new purchases/payments entries.
# Order matters so that new shared strings land at the expected indices
# (21 = "Phone Case", 22 = "New Clothes", 23 = "Check 7/14/16").
$ws.Range("B9").Value = "Phone Case"
$ws.Range("F7").Value = "New Clothes"
$ws.Range("H15").Value = "Check 7/14/16"

$ws.Range("B10").Value = 69.24
$ws.Range("F6").Value = "Purchases"
$ws.Range("F8").Value = 66.82
$ws.Range("H16").Value = 252.88

# Recalculate so dependent formulas (A2:H2, A4:H4, A5:H5) update.
$wb.Application.Calculate()

# Update the selected cell shown in the sheet view.
$ws.Range("K13").Select() | Out-Null
